$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.871.48"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.902.32"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8031"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.02"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.908.79"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7330"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.153"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.51"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "29.877.51"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.85"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.830"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.66"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007687"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "2.155.26"
$ws.Range("E22").Value = "  +0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.874"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.29"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.169"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1408"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +7.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.78"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.008"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.360"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.513"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.270"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05546"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.048"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.250"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7265"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.42%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01917"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.789"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4383"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.78%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.039"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.73"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8351"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.854"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.38"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.523"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.688"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").Value = "2.062.67"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "976.76"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.03"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.54%  "
